$d = $word.ActiveDocument

# Paragraph 12 is the "Kính gửi:  Ban Giám hiệu Trường Đại học Xây dựng<tab>"
# salutation line. Replace the trailing tab character with " Hà Nội" so the
# line reads "...Trường Đại học Xây dựng Hà Nội".
$p = $d.Paragraphs.Item(12)
$full = $p.Range
$tabRange = $d.Range($full.End - 2, $full.End - 1)
$tabRange.Text = " Hà Nội"
